$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 43
$ws.Cells.Item(43, 8).Value = 5462.375  # H43: 5188.778 -> 5462.375
$ws.Cells.Item(43, 10).Value = 5933  # J43: 5199.75 -> 5933
$ws.Cells.Item(43, 12).Value = 5933  # L43: 5199.75 -> 5933
$ws.Cells.Item(43, 14).Value = -6071  # N43: -5337.75 -> -6071
# Row 74
$ws.Cells.Item(74, 8).Value = 72717496  # H74: 72717460 -> 72717496
$ws.Cells.Item(74, 10).Value = 7250  # J74: 7100 -> 7250
$ws.Cells.Item(74, 12).Value = 7250  # L74: 7100 -> 7250
$ws.Cells.Item(74, 14).Value = -9122  # N74: -8972 -> -9122
# Row 77
$ws.Cells.Item(77, 8).Value = 72717496  # H77: 72717460 -> 72717496
$ws.Cells.Item(77, 10).Value = 7250  # J77: 7100 -> 7250
$ws.Cells.Item(77, 12).Value = 36250  # L77: 35500 -> 36250
$ws.Cells.Item(77, 14).Value = -45610  # N77: -44860 -> -45610
# Row 86
$ws.Cells.Item(86, 8).Value = 9306.066000000001  # H86: 10007.077 -> 9306.066000000001
$ws.Cells.Item(86, 10).Value = 5236.75  # J86: 5724 -> 5236.75
$ws.Cells.Item(86, 12).Value = 5236.75  # L86: 5724 -> 5236.75
$ws.Cells.Item(86, 14).Value = -7482.75  # N86: -7970 -> -7482.75
# Row 89
$ws.Cells.Item(89, 8).Value = 9306.066000000001  # H89: 10007.077 -> 9306.066000000001
$ws.Cells.Item(89, 10).Value = 5236.75  # J89: 5724 -> 5236.75
$ws.Cells.Item(89, 12).Value = 26183.75  # L89: 28620 -> 26183.75
$ws.Cells.Item(89, 14).Value = -37415.75  # N89: -39852 -> -37415.75
# Row 100
$ws.Cells.Item(100, 8).Value = 4833.1665  # H100: 5999.5 -> 4833.1665
$ws.Cells.Item(100, 9).Value = 5200.2  # I100: 5999.5 -> 5200.2
$ws.Cells.Item(100, 10).Value = 2998  # J100: 0 -> 2998
$ws.Cells.Item(100, 11).Value = 5200.2  # K100: 5999.5 -> 5200.2
$ws.Cells.Item(100, 12).Value = 2998  # L100: 0 -> 2998
$ws.Cells.Item(100, 13).Value = -4659.2  # M100: -5458.5 -> -4659.2
$ws.Cells.Item(100, 14).Value = -4080  # N100: None -> -4080

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 102
$ws.Cells.Item(102, 8).Value = 6251395.5  # H102: 6668042 -> 6251395.5
$ws.Cells.Item(102, 10).Value = 1431.3334  # J102: 1299.5 -> 1431.3334
$ws.Cells.Item(102, 12).Value = 1431.3334  # L102: 1299.5 -> 1431.3334
$ws.Cells.Item(102, 14).Value = -4675.3334  # N102: -4543.5 -> -4675.3334

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Cells.Item(99, 8).Value = 2017.3158  # H99: 2004.5385 -> 2017.3158
$ws.Cells.Item(99, 9).Value = 1944.5834  # I99: 1927.56 -> 1944.5834
$ws.Cells.Item(99, 11).Value = 1944.5834  # K99: 1927.56 -> 1944.5834
$ws.Cells.Item(99, 13).Value = -446.5834  # M99: -429.5599999999999 -> -446.5834
# Row 107
$ws.Cells.Item(107, 8).Value = 169794.83  # H107: 253597.25 -> 169794.83
$ws.Cells.Item(107, 10).Value = 203313.8  # J107: 337396.34 -> 203313.8
$ws.Cells.Item(107, 12).Value = 203313.8  # L107: 337396.34 -> 203313.8
$ws.Cells.Item(107, 14).Value = -207153.8  # N107: -341236.34 -> -207153.8

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Cells.Item(16, 8).Value = 2719649.8  # H16: 3625196.2 -> 2719649.8
$ws.Cells.Item(16, 9).Value = 3625195.2  # I16: 5436288 -> 3625195.2
$ws.Cells.Item(16, 11).Value = 3625195.2  # K16: 5436288 -> 3625195.2
$ws.Cells.Item(16, 13).Value = -3624908.2  # M16: -5436001 -> -3624908.2
# Row 31
$ws.Cells.Item(31, 8).Value = 4508.3447  # H31: 4543.6206 -> 4508.3447
$ws.Cells.Item(31, 9).Value = 3147  # I31: 3063.8823 -> 3147
$ws.Cells.Item(31, 10).Value = 6183.846  # J31: 6639.9165 -> 6183.846
$ws.Cells.Item(31, 11).Value = 3147  # K31: 3063.8823 -> 3147
$ws.Cells.Item(31, 12).Value = 6183.846  # L31: 6639.9165 -> 6183.846
$ws.Cells.Item(31, 13).Value = -2852  # M31: -2768.8823 -> -2852
$ws.Cells.Item(31, 14).Value = -6773.846  # N31: -7229.9165 -> -6773.846
# Row 34
$ws.Cells.Item(34, 8).Value = 4508.3447  # H34: 4543.6206 -> 4508.3447
$ws.Cells.Item(34, 9).Value = 3147  # I34: 3063.8823 -> 3147
$ws.Cells.Item(34, 10).Value = 6183.846  # J34: 6639.9165 -> 6183.846
$ws.Cells.Item(34, 11).Value = 3147  # K34: 3063.8823 -> 3147
$ws.Cells.Item(34, 12).Value = 6183.846  # L34: 6639.9165 -> 6183.846
$ws.Cells.Item(34, 13).Value = -2945  # M34: -2861.8823 -> -2945
$ws.Cells.Item(34, 14).Value = -6587.846  # N34: -7043.9165 -> -6587.846
# Row 58
$ws.Cells.Item(58, 8).Value = 26328590  # H58: 27791006 -> 26328590
$ws.Cells.Item(58, 9).Value = 26328590  # I58: 29425184 -> 26328590
$ws.Cells.Item(58, 10).Value = 0  # J58: 9999 -> 0
$ws.Cells.Item(58, 11).Value = 26328590  # K58: 29425184 -> 26328590
$ws.Cells.Item(58, 12).Value = 0  # L58: 9999 -> 0
$ws.Cells.Item(58, 13).ClearContents()  # M58: remove (was -29424981)
$ws.Cells.Item(58, 14).Value = -26328387  # N58: -10405 -> -26328387
# Row 68
$ws.Cells.Item(68, 8).Value = 87875  # H68: 90099.8 -> 87875
$ws.Cells.Item(68, 10).Value = 99000  # J68: 98999.78 -> 99000
$ws.Cells.Item(68, 12).Value = 99000  # L68: 98999.78 -> 99000
$ws.Cells.Item(68, 14).Value = -100498  # N68: -100497.78 -> -100498
# Row 71
$ws.Cells.Item(71, 8).Value = 87875  # H71: 90099.8 -> 87875
$ws.Cells.Item(71, 10).Value = 99000  # J71: 98999.78 -> 99000
$ws.Cells.Item(71, 12).Value = 297000  # L71: 296999.34 -> 297000
$ws.Cells.Item(71, 14).Value = -304488  # N71: -304487.34 -> -304488
# Row 113
$ws.Cells.Item(113, 8).Value = 2719649.8  # H113: 3625196.2 -> 2719649.8
$ws.Cells.Item(113, 9).Value = 3625195.2  # I113: 5436288 -> 3625195.2
$ws.Cells.Item(113, 11).Value = 3625195.2  # K113: 5436288 -> 3625195.2
$ws.Cells.Item(113, 13).Value = -3623025.2  # M113: -5434118 -> -3623025.2
# Row 136
$ws.Cells.Item(136, 8).Value = 26328590  # H136: 27791006 -> 26328590
$ws.Cells.Item(136, 9).Value = 26328590  # I136: 29425184 -> 26328590
$ws.Cells.Item(136, 10).Value = 0  # J136: 9999 -> 0
$ws.Cells.Item(136, 11).Value = 78985770  # K136: 88275552 -> 78985770
$ws.Cells.Item(136, 12).Value = 0  # L136: 29997 -> 0
$ws.Cells.Item(136, 13).ClearContents()  # M136: remove (was -88273002)
$ws.Cells.Item(136, 14).Value = -78983220  # N136: -35097 -> -78983220

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Cells.Item(5, 8).Value = 112200  # H5: 101035.8 -> 112200
$ws.Cells.Item(5, 9).Value = 200453.2  # I5: 167137.33 -> 200453.2
$ws.Cells.Item(5, 11).Value = 601359.6000000001  # K5: 501411.99 -> 601359.6000000001
$ws.Cells.Item(5, 13).Value = -601247.6000000001  # M5: -501299.99 -> -601247.6000000001
# Row 87
$ws.Cells.Item(87, 8).Value = 922.8  # H87: 953.75 -> 922.8
$ws.Cells.Item(87, 9).Value = 899.5  # I87: 933 -> 899.5
$ws.Cells.Item(87, 11).Value = 2698.5  # K87: 2799 -> 2698.5
$ws.Cells.Item(87, 13).Value = -1450.5  # M87: -1551 -> -1450.5
# Row 90
$ws.Cells.Item(90, 8).Value = 922.8  # H90: 953.75 -> 922.8
$ws.Cells.Item(90, 9).Value = 899.5  # I90: 933 -> 899.5
$ws.Cells.Item(90, 11).Value = 8095.5  # K90: 8397 -> 8095.5
$ws.Cells.Item(90, 13).Value = -1855.5  # M90: -2157 -> -1855.5
# Row 109
$ws.Cells.Item(109, 8).Value = 2521  # H109: 3100.8333 -> 2521
$ws.Cells.Item(109, 10).Value = 0  # J109: 6000 -> 0
$ws.Cells.Item(109, 12).Value = 0  # L109: 18000 -> 0
$ws.Cells.Item(109, 14).ClearContents()  # N109: remove (was -20080)
# Row 113
$ws.Cells.Item(113, 8).Value = 111681.555  # H113: 111681.664 -> 111681.555
$ws.Cells.Item(113, 10).Value = 836.6  # J113: 836.8 -> 836.6
$ws.Cells.Item(113, 12).Value = 2509.8  # L113: 2510.4 -> 2509.8
$ws.Cells.Item(113, 14).Value = -6849.8  # N113: -6850.4 -> -6849.8
# Row 119
$ws.Cells.Item(119, 8).Value = 5305.4  # H119: 1850 -> 5305.4
$ws.Cells.Item(119, 9).Value = 1632  # I119: 1850 -> 1632
$ws.Cells.Item(119, 10).Value = 19999  # J119: 0 -> 19999
$ws.Cells.Item(119, 11).Value = 4896  # K119: 5550 -> 4896
$ws.Cells.Item(119, 12).Value = 59997  # L119: 0 -> 59997
$ws.Cells.Item(119, 13).Value = -58  # M119: -712 -> -58
$ws.Cells.Item(119, 14).Value = -69673  # N119: None -> -69673
# Row 129
$ws.Cells.Item(129, 8).Value = 3489.52  # H129: 3823.5 -> 3489.52
$ws.Cells.Item(129, 9).Value = 396.66666  # I129: 400 -> 396.66666
$ws.Cells.Item(129, 10).Value = 3911.2727  # J129: 4165.85 -> 3911.2727
$ws.Cells.Item(129, 11).Value = 1189.99998  # K129: 1200 -> 1189.99998
$ws.Cells.Item(129, 12).Value = 11733.8181  # L129: 12497.55 -> 11733.8181
$ws.Cells.Item(129, 13).Value = 3810.00002  # M129: 3800 -> 3810.00002
$ws.Cells.Item(129, 14).Value = -21733.8181  # N129: -22497.55 -> -21733.8181
# Row 135
$ws.Cells.Item(135, 8).Value = 112200  # H135: 101035.8 -> 112200
$ws.Cells.Item(135, 9).Value = 200453.2  # I135: 167137.33 -> 200453.2
$ws.Cells.Item(135, 11).Value = 1804078.8  # K135: 1504235.97 -> 1804078.8
$ws.Cells.Item(135, 13).Value = -1801543.8  # M135: -1501700.97 -> -1801543.8

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Cells.Item(122, 8).Value = 5520.3447  # H122: 5630.8623 -> 5520.3447
$ws.Cells.Item(122, 9).Value = 3534.5217  # I122: 3673.8696 -> 3534.5217
$ws.Cells.Item(122, 11).Value = 10603.5651  # K122: 11021.6088 -> 10603.5651
$ws.Cells.Item(122, 13).Value = -8153.5651  # M122: -8571.6088 -> -8153.5651
# Row 123
$ws.Cells.Item(123, 8).Value = 0  # H123: 96990 -> 0
$ws.Cells.Item(123, 10).Value = 0  # J123: 96990 -> 0
$ws.Cells.Item(123, 12).ClearContents()  # L123: remove (was 96990)
$ws.Cells.Item(123, 14).Value = 0  # N123: -101890 -> 0
# Row 132
$ws.Cells.Item(132, 8).Value = 3682552.5  # H132: 3214599.2 -> 3682552.5
$ws.Cells.Item(132, 9).Value = 3912384.2  # I132: 3482052 -> 3912384.2
$ws.Cells.Item(132, 10).Value = 5247  # J132: 5164.6665 -> 5247
$ws.Cells.Item(132, 11).Value = 11737152.6  # K132: 10446156 -> 11737152.6
$ws.Cells.Item(132, 12).Value = 15741  # L132: 15493.9995 -> 15741
$ws.Cells.Item(132, 13).Value = -11734622.6  # M132: -10443626 -> -11734622.6
$ws.Cells.Item(132, 14).Value = -20801  # N132: -20553.9995 -> -20801

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Cells.Item(55, 8).Value = 556.7857  # H55: 552.7857 -> 556.7857
$ws.Cells.Item(55, 9).Value = 356  # I55: 309.2 -> 356
$ws.Cells.Item(55, 10).Value = 637.1  # J55: 688.1111 -> 637.1
$ws.Cells.Item(55, 11).Value = 356  # K55: 309.2 -> 356
$ws.Cells.Item(55, 12).Value = 637.1  # L55: 688.1111 -> 637.1
$ws.Cells.Item(55, 13).Value = -183  # M55: -136.2 -> -183
$ws.Cells.Item(55, 14).Value = -983.1  # N55: -1034.1111 -> -983.1
# Row 93
$ws.Cells.Item(93, 8).Value = 1675.9286  # H93: 1775 -> 1675.9286
$ws.Cells.Item(93, 9).Value = 1164  # I93: 1230.6666 -> 1164
$ws.Cells.Item(93, 10).Value = 2597.4  # J93: 2999.75 -> 2597.4
$ws.Cells.Item(93, 11).Value = 1164  # K93: 1230.6666 -> 1164
$ws.Cells.Item(93, 12).Value = 2597.4  # L93: 2999.75 -> 2597.4
$ws.Cells.Item(93, 13).Value = 84  # M93: 17.33339999999998 -> 84
$ws.Cells.Item(93, 14).Value = -5093.4  # N93: -5495.75 -> -5093.4
# Row 132
$ws.Cells.Item(132, 8).Value = 7357873  # H132: 7357873.5 -> 7357873
$ws.Cells.Item(132, 9).Value = 8626009  # I132: 8338507.5 -> 8626009
$ws.Cells.Item(132, 10).Value = 2687.6  # J132: 3122 -> 2687.6
$ws.Cells.Item(132, 11).Value = 25878027  # K132: 25015522.5 -> 25878027
$ws.Cells.Item(132, 12).Value = 8062.799999999999  # L132: 9366 -> 8062.799999999999
$ws.Cells.Item(132, 13).Value = -25875497  # M132: -25012992.5 -> -25875497
$ws.Cells.Item(132, 14).Value = -13122.8  # N132: -14426 -> -13122.8
# Row 136
$ws.Cells.Item(136, 8).Value = 1485.3158  # H136: 1541.4722 -> 1485.3158
$ws.Cells.Item(136, 9).Value = 1232.3928  # I136: 1316.5927 -> 1232.3928
$ws.Cells.Item(136, 10).Value = 2193.5  # J136: 2216.111 -> 2193.5
$ws.Cells.Item(136, 11).Value = 3697.1784  # K136: 3949.7781 -> 3697.1784
$ws.Cells.Item(136, 12).Value = 6580.5  # L136: 6648.333 -> 6580.5
$ws.Cells.Item(136, 13).Value = -1147.1784  # M136: -1399.7781 -> -1147.1784
$ws.Cells.Item(136, 14).Value = -11680.5  # N136: -11748.333 -> -11680.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Cells.Item(107, 8).Value = 983.5  # H107: 991.1667 -> 983.5
$ws.Cells.Item(107, 9).Value = 788.0909  # I107: 799.4 -> 788.0909
$ws.Cells.Item(107, 10).Value = 1700  # J107: 1950 -> 1700
$ws.Cells.Item(107, 11).Value = 2364.2727  # K107: 2398.2 -> 2364.2727
$ws.Cells.Item(107, 12).Value = 5100  # L107: 5850 -> 5100
$ws.Cells.Item(107, 13).Value = -444.2727  # M107: -478.1999999999998 -> -444.2727
$ws.Cells.Item(107, 14).Value = -8940  # N107: -9690 -> -8940
